# Insert two new data rows (796-797) into the "Cebolla" sheet, pushing the
# existing rows 796..892 down to 798..894. This models a new weekly price
# entry being added near the top of the existing block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 796, shifting everything from old row 796 onward
# down by two (old 796 -> 798, ..., old 892 -> 894).
$ws.Rows("796:797").Insert()

# ---- New row 796 ----
$ws.Range("A796").Value = 11
$ws.Range("B796").Value = "Vega Monumental Concepción"
$ws.Range("C796").Value = "Bíobío"
$ws.Range("D796").Value = 45142
$ws.Range("E796").Value = 8
$ws.Range("F796").Value = 100112004
$ws.Range("G796").Value = "Cebolla"
$ws.Range("H796").Value = "Morada(o)"
$ws.Range("I796").Value = "1a nueva(o)"
$ws.Range("J796").Value = 80
$ws.Range("K796").Value = 12000
$ws.Range("L796").Value = 12000
$ws.Range("M796").Value = 12000
$ws.Range("N796").Value = "`$/malla 18 kilos"
$ws.Range("O796").Value = "Región de Arica y Parinacota"
$ws.Range("P796").Value = 667
$ws.Range("Q796").Value = 18
$ws.Range("R796").Value = "Hortaliza"

# ---- New row 797 ----
$ws.Range("A797").Value = 11
$ws.Range("B797").Value = "Vega Monumental Concepción"
$ws.Range("C797").Value = "Bíobío"
$ws.Range("D797").Value = 45142
$ws.Range("E797").Value = 8
$ws.Range("F797").Value = 100112004
$ws.Range("G797").Value = "Cebolla"
$ws.Range("H797").Value = "Sin especificar"
$ws.Range("I797").Value = "1a (guarda)"
$ws.Range("J797").Value = 250
$ws.Range("K797").Value = 10000
$ws.Range("L797").Value = 11000
$ws.Range("M797").Value = 10400
$ws.Range("N797").Value = "`$/malla 18 kilos"
$ws.Range("O797").Value = "Región de O'Higgins"
$ws.Range("P797").Value = 578
$ws.Range("Q797").Value = 18
$ws.Range("R797").Value = "Hortaliza"
